$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("Gesamtübersicht").Name = "Complete Overview"
$wb.Worksheets.Item("Monatsübersicht").Name = "Month Overview"
